$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "1. Clic en el boton ID DEAL`n2. Ingresar número ID DEAL válido`n3. clic en Consultar cliente"
$ws.Range("F4").Value = "1. Clic en el boton OPCIONES`n2. Clic en opción `"Configuración WiFi`"`n3. Seleccionar el campo `"Nombre de red`"`n4. Digitar nuevo nombre de red`n5. Clic en el select CANAL`n6. Selección aleatoria de canal`n7. Clic en select ANCHO BANDA CANAL`n8. Selección aleatoria de ancho de banda`n9. Marcar checkbox 'Unsecured'`n10. Clic en botón ENVIAR y esperar progress`n11. Cerrar modal de Configuración WiFi"
$ws.Range("B7").Value = "funcion UPnP(opcion click boton cancelar)"
$ws.Range("B8").Value = "funcion DMZ(opcion click boton cancelar)"
$ws.Range("B9").Value = "ipv4 port Mapping(opcion click boton cancelar)"
$ws.Range("B10").Value = "reserva DHCP(opcion click boton cancelar)"

$ws.Range("A11").Select()
